$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9:D9").Font.Color = 0
$ws.Range("C10").Font.Color = 0

$ws.Range("A9:D9").RowHeight = 47.25
$ws.Range("A10:D10").RowHeight = 47.25
